$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.169591333333333
$ws.Range("H2").Value = 3.508774
$ws.Range("I2").Value = 0.05593990076588554
$ws.Range("J2").Value = 0.05593990076588554
$ws.Range("M2").Value = 106.5570066666667
$ws.Range("N2").Value = 319.67102
$ws.Range("O2").Value = 0.956549115742331
$ws.Range("P2").Value = 0.956549115742331
$ws.Range("Q2").Value = 124.6281515032755
$ws.Range("R2").Value = 1121.65336352948
$ws.Range("S2").Value = 0.05350926261232156
$ws.Range("T2").Value = 0.05350926261232156
$ws.Range("G3").Value = 1.169591333333333
$ws.Range("H3").Value = 3.508774
$ws.Range("I3").Value = 0.05593990076588554
$ws.Range("J3").Value = 0.05593990076588554
$ws.Range("M3").Value = 0.249583
$ws.Range("N3").Value = 0.7487489999999999
$ws.Range("O3").Value = 0.002240475830004717
$ws.Range("P3").Value = 0.002240475830004717
$ws.Range("Q3").Value = 0.2919101137473333
$ws.Range("R3").Value = 2.627191023726
$ws.Range("S3").Value = 0.0001253319955988289
$ws.Range("T3").Value = 0.0001253319955988289
$ws.Range("G4").Value = 1.169591333333333
$ws.Range("H4").Value = 3.508774
$ws.Range("I4").Value = 0.05593990076588554
$ws.Range("J4").Value = 0.05593990076588554
$ws.Range("M4").Value = 3.632925333333334
$ws.Range("N4").Value = 10.898776
$ws.Range("O4").Value = 0.03261232296087941
$ws.Range("P4").Value = 0.03261232296087941
$ws.Range("Q4").Value = 4.249037984513778
$ws.Range("R4").Value = 38.241341860624
$ws.Range("S4").Value = 0.001824330110176604
$ws.Range("T4").Value = 0.001824330110176604
$ws.Range("G5").Value = 1.169591333333333
$ws.Range("H5").Value = 3.508774
$ws.Range("I5").Value = 0.05593990076588554
$ws.Range("J5").Value = 0.05593990076588554
$ws.Range("M5").Value = 0.7472223333333332
$ws.Range("N5").Value = 2.241667
$ws.Range("O5").Value = 0.006707722791508481
$ws.Range("P5").Value = 0.006707722791508481
$ws.Range("Q5").Value = 0.8739447651397776
$ws.Range("R5").Value = 7.865502886257999
$ws.Range("S5").Value = 0.0003752293473220532
$ws.Range("T5").Value = 0.0003752293473220532
$ws.Range("G6").Value = 1.169591333333333
$ws.Range("H6").Value = 3.508774
$ws.Range("I6").Value = 0.05593990076588554
$ws.Range("J6").Value = 0.05593990076588554
$ws.Range("M6").Value = 0.2105813333333333
$ws.Range("N6").Value = 0.631744
$ws.Range("O6").Value = 0.001890362675276361
$ws.Range("P6").Value = 0.001890362675276361
$ws.Range("Q6").Value = 0.2462941024284444
$ws.Range("R6").Value = 2.216646921856
$ws.Range("S6").Value = 0.0001057467004664935
$ws.Range("T6").Value = 0.0001057467004664935
$ws.Range("G7").Value = 5.220282666666667
$ws.Range("H7").Value = 15.660848
$ws.Range("I7").Value = 0.2496787433529823
$ws.Range("J7").Value = 0.2496787433529823
$ws.Range("M7").Value = 106.5570066666667
$ws.Range("N7").Value = 319.67102
$ws.Range("O7").Value = 0.956549115742331
$ws.Range("P7").Value = 0.956549115742331
$ws.Range("Q7").Value = 556.2576949138845
$ws.Range("R7").Value = 5006.31925422496
$ws.Range("S7").Value = 0.2388299811739516
$ws.Range("T7").Value = 0.2388299811739516
$ws.Range("G8").Value = 5.220282666666667
$ws.Range("H8").Value = 15.660848
$ws.Range("I8").Value = 0.2496787433529823
$ws.Range("J8").Value = 0.2496787433529823
$ws.Range("M8").Value = 0.249583
$ws.Range("N8").Value = 0.7487489999999999
$ws.Range("O8").Value = 0.002240475830004717
$ws.Range("P8").Value = 0.002240475830004717
$ws.Range("Q8").Value = 1.302893808794667
$ws.Range("R8").Value = 11.726044279152
$ws.Range("S8").Value = 0.0005593991897483078
$ws.Range("T8").Value = 0.0005593991897483076
$ws.Range("G9").Value = 5.220282666666667
$ws.Range("H9").Value = 15.660848
$ws.Range("I9").Value = 0.2496787433529823
$ws.Range("J9").Value = 0.2496787433529823
$ws.Range("M9").Value = 3.632925333333334
$ws.Range("N9").Value = 10.898776
$ws.Range("O9").Value = 0.03261232296087941
$ws.Range("P9").Value = 0.03261232296087941
$ws.Range("Q9").Value = 18.96489714689423
$ws.Range("R9").Value = 170.684074322048
$ws.Range("S9").Value = 0.008142603814693982
$ws.Range("T9").Value = 0.00814260381469398
$ws.Range("G10").Value = 5.220282666666667
$ws.Range("H10").Value = 15.660848
$ws.Range("I10").Value = 0.2496787433529823
$ws.Range("J10").Value = 0.2496787433529823
$ws.Range("M10").Value = 0.7472223333333332
$ws.Range("N10").Value = 2.241667
$ws.Range("O10").Value = 0.006707722791508481
$ws.Range("P10").Value = 0.006707722791508481
$ws.Range("Q10").Value = 3.900711794846222
$ws.Range("R10").Value = 35.10640615361599
$ws.Range("S10").Value = 0.001674775797343996
$ws.Range("T10").Value = 0.001674775797343996
$ws.Range("G11").Value = 5.220282666666667
$ws.Range("H11").Value = 15.660848
$ws.Range("I11").Value = 0.2496787433529823
$ws.Range("J11").Value = 0.2496787433529823
$ws.Range("M11").Value = 0.2105813333333333
$ws.Range("N11").Value = 0.631744
$ws.Range("O11").Value = 0.001890362675276361
$ws.Range("P11").Value = 0.001890362675276361
$ws.Range("Q11").Value = 1.099294084323555
$ws.Range("R11").Value = 9.893646758911999
$ws.Range("S11").Value = 0.0004719833772443836
$ws.Range("T11").Value = 0.0004719833772443835
$ws.Range("G12").Value = 4.885583666666666
$ws.Range("H12").Value = 14.656751
$ws.Range("I12").Value = 0.2336705631341014
$ws.Range("J12").Value = 0.2336705631341014
$ws.Range("M12").Value = 106.5570066666667
$ws.Range("N12").Value = 319.67102
$ws.Range("O12").Value = 0.956549115742331
$ws.Range("P12").Value = 0.956549115742331
$ws.Range("Q12").Value = 520.5931713395578
$ws.Range("R12").Value = 4685.33854205602
$ws.Range("S12").Value = 0.2235173705409373
$ws.Range("T12").Value = 0.2235173705409373
$ws.Range("G13").Value = 4.885583666666666
$ws.Range("H13").Value = 14.656751
$ws.Range("I13").Value = 0.2336705631341014
$ws.Range("J13").Value = 0.2336705631341014
$ws.Range("M13").Value = 0.249583
$ws.Range("N13").Value = 0.7487489999999999
$ws.Range("O13").Value = 0.002240475830004717
$ws.Range("P13").Value = 0.002240475830004717
$ws.Range("Q13").Value = 1.219358628277667
$ws.Range("R13").Value = 10.974227654499
$ws.Range("S13").Value = 0.0005235332488855456
$ws.Range("T13").Value = 0.0005235332488855455
$ws.Range("G14").Value = 4.885583666666666
$ws.Range("H14").Value = 14.656751
$ws.Range("I14").Value = 0.2336705631341014
$ws.Range("J14").Value = 0.2336705631341014
$ws.Range("M14").Value = 3.632925333333334
$ws.Range("N14").Value = 10.898776
$ws.Range("O14").Value = 0.03261232296087941
$ws.Range("P14").Value = 0.03261232296087941
$ws.Range("Q14").Value = 17.74896067075289
$ws.Range("R14").Value = 159.740646036776
$ws.Range("S14").Value = 0.007620539871379878
$ws.Range("T14").Value = 0.007620539871379878
$ws.Range("G15").Value = 4.885583666666666
$ws.Range("H15").Value = 14.656751
$ws.Range("I15").Value = 0.2336705631341014
$ws.Range("J15").Value = 0.2336705631341014
$ws.Range("M15").Value = 0.7472223333333332
$ws.Range("N15").Value = 2.241667
$ws.Range("O15").Value = 0.006707722791508481
$ws.Range("P15").Value = 0.006707722791508481
$ws.Range("Q15").Value = 3.650617227101888
$ws.Range("R15").Value = 32.85555504391699
$ws.Range("S15").Value = 0.001567397362039234
$ws.Range("T15").Value = 0.001567397362039234
$ws.Range("G16").Value = 4.885583666666666
$ws.Range("H16").Value = 14.656751
$ws.Range("I16").Value = 0.2336705631341014
$ws.Range("J16").Value = 0.2336705631341014
$ws.Range("M16").Value = 0.2105813333333333
$ws.Range("N16").Value = 0.631744
$ws.Range("O16").Value = 0.001890362675276361
$ws.Range("P16").Value = 0.001890362675276361
$ws.Range("Q16").Value = 1.028812722638222
$ws.Range("R16").Value = 9.259314503744
$ws.Range("S16").Value = 0.0004417221108595138
$ws.Range("T16").Value = 0.0004417221108595138
$ws.Range("G17").Value = 2.891468666666666
$ws.Range("H17").Value = 8.674406
$ws.Range("I17").Value = 0.1382948604962879
$ws.Range("J17").Value = 0.1382948604962879
$ws.Range("M17").Value = 106.5570066666667
$ws.Range("N17").Value = 319.67102
$ws.Range("O17").Value = 0.956549115742331
$ws.Range("P17").Value = 0.956549115742331
$ws.Range("Q17").Value = 308.1062459904578
$ws.Range("R17").Value = 2772.95621391412
$ws.Range("S17").Value = 0.1322858265194332
$ws.Range("T17").Value = 0.1322858265194332
$ws.Range("G18").Value = 2.891468666666666
$ws.Range("H18").Value = 8.674406
$ws.Range("I18").Value = 0.1382948604962879
$ws.Range("J18").Value = 0.1382948604962879
$ws.Range("M18").Value = 0.249583
$ws.Range("N18").Value = 0.7487489999999999
$ws.Range("O18").Value = 0.002240475830004717
$ws.Range("P18").Value = 0.002240475830004717
$ws.Range("Q18").Value = 0.7216614242326665
$ws.Range("R18").Value = 6.494952818093998
$ws.Range("S18").Value = 0.0003098462923558073
$ws.Range("T18").Value = 0.0003098462923558072
$ws.Range("G19").Value = 2.891468666666666
$ws.Range("H19").Value = 8.674406
$ws.Range("I19").Value = 0.1382948604962879
$ws.Range("J19").Value = 0.1382948604962879
$ws.Range("M19").Value = 3.632925333333334
$ws.Range("N19").Value = 10.898776
$ws.Range("O19").Value = 0.03261232296087941
$ws.Range("P19").Value = 0.03261232296087941
$ws.Range("Q19").Value = 10.50448976967289
$ws.Range("R19").Value = 94.54040792705601
$ws.Range("S19").Value = 0.004510116654334705
$ws.Range("T19").Value = 0.004510116654334704
$ws.Range("G20").Value = 2.891468666666666
$ws.Range("H20").Value = 8.674406
$ws.Range("I20").Value = 0.1382948604962879
$ws.Range("J20").Value = 0.1382948604962879
$ws.Range("M20").Value = 0.7472223333333332
$ws.Range("N20").Value = 2.241667
$ws.Range("O20").Value = 0.006707722791508481
$ws.Range("P20").Value = 0.006707722791508481
$ws.Range("Q20").Value = 2.160569963866888
$ws.Range("R20").Value = 19.445129674802
$ws.Range("S20").Value = 0.0009276435876994364
$ws.Range("T20").Value = 0.0009276435876994362
$ws.Range("G21").Value = 2.891468666666666
$ws.Range("H21").Value = 8.674406
$ws.Range("I21").Value = 0.1382948604962879
$ws.Range("J21").Value = 0.1382948604962879
$ws.Range("M21").Value = 0.2105813333333333
$ws.Range("N21").Value = 0.631744
$ws.Range("O21").Value = 0.001890362675276361
$ws.Range("P21").Value = 0.001890362675276361
$ws.Range("Q21").Value = 0.6088893271182221
$ws.Range("R21").Value = 5.480003944063999
$ws.Range("S21").Value = 0.000261427442464734
$ws.Range("T21").Value = 0.0002614274424647339
$ws.Range("G22").Value = 6.741071666666667
$ws.Range("H22").Value = 20.223215
$ws.Range("I22").Value = 0.3224159322507428
$ws.Range("J22").Value = 0.3224159322507428
$ws.Range("M22").Value = 106.5570066666667
$ws.Range("N22").Value = 319.67102
$ws.Range("O22").Value = 0.956549115742331
$ws.Range("P22").Value = 0.956549115742331
$ws.Range("Q22").Value = 718.3084185254778
$ws.Range("R22").Value = 6464.7757667293
$ws.Range("S22").Value = 0.3084066748956873
$ws.Range("T22").Value = 0.3084066748956873
$ws.Range("G23").Value = 6.741071666666667
$ws.Range("H23").Value = 20.223215
$ws.Range("I23").Value = 0.3224159322507428
$ws.Range("J23").Value = 0.3224159322507428
$ws.Range("M23").Value = 0.249583
$ws.Range("N23").Value = 0.7487489999999999
$ws.Range("O23").Value = 0.002240475830004717
$ws.Range("P23").Value = 0.002240475830004717
$ws.Range("Q23").Value = 1.682456889781667
$ws.Range("R23").Value = 15.142112008035
$ws.Range("S23").Value = 0.0007223651034162277
$ws.Range("T23").Value = 0.0007223651034162276
$ws.Range("G24").Value = 6.741071666666667
$ws.Range("H24").Value = 20.223215
$ws.Range("I24").Value = 0.3224159322507428
$ws.Range("J24").Value = 0.3224159322507428
$ws.Range("M24").Value = 3.632925333333334
$ws.Range("N24").Value = 10.898776
$ws.Range("O24").Value = 0.03261232296087941
$ws.Range("P24").Value = 0.03261232296087941
$ws.Range("Q24").Value = 24.48981003164889
$ws.Range("R24").Value = 220.40829028484
$ws.Range("S24").Value = 0.01051473251029424
$ws.Range("T24").Value = 0.01051473251029424
$ws.Range("G25").Value = 6.741071666666667
$ws.Range("H25").Value = 20.223215
$ws.Range("I25").Value = 0.3224159322507428
$ws.Range("J25").Value = 0.3224159322507428
$ws.Range("M25").Value = 0.7472223333333332
$ws.Range("N25").Value = 2.241667
$ws.Range("O25").Value = 0.006707722791508481
$ws.Range("P25").Value = 0.006707722791508481
$ws.Range("Q25").Value = 5.037079299933888
$ws.Range("R25").Value = 45.33371369940499
$ws.Range("S25").Value = 0.002162676697103762
$ws.Range("T25").Value = 0.002162676697103762
$ws.Range("G26").Value = 6.741071666666667
$ws.Range("H26").Value = 20.223215
$ws.Range("I26").Value = 0.3224159322507428
$ws.Range("J26").Value = 0.3224159322507428
$ws.Range("M26").Value = 0.2105813333333333
$ws.Range("N26").Value = 0.631744
$ws.Range("O26").Value = 0.001890362675276361
$ws.Range("P26").Value = 0.001890362675276361
$ws.Range("Q26").Value = 1.419543859662222
$ws.Range("R26").Value = 12.77589473696
$ws.Range("S26").Value = 0.0006094830442412362
$ws.Range("T26").Value = 0.0006094830442412362

Write-Output "Done applying changes."